# ---------------------------------------------------------------------------
# Applies the "Updated App with Instructions as first tab. Fixed issues with
# multiple curves being returned" commit to AdminData.xlsx.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) DatePeriods (sheet1): trim the reporting-period lookup table down to a
#    single quarter (the old "2017 Q3" row, which actually carries the
#    Oct-Dec date range) plus the existing trailing blank entry row.
# ---------------------------------------------------------------------------
$wsDates = $wb.Worksheets.Item("DatePeriods")
$wsDates.Rows("5:13").Delete()
$wsDates.Rows("2:3").Delete()

# ---------------------------------------------------------------------------
# 2) OrgSetUp (sheet2): keep only Area1/Sub1, Area1/Sub2 and Area2/Sub6.
# ---------------------------------------------------------------------------
$wsOrg = $wb.Worksheets.Item("OrgSetUp")
$wsOrg.Rows("8:16").Delete()
$wsOrg.Rows("4:6").Delete()

# ---------------------------------------------------------------------------
# 3) RespCurveLookup (sheet4): the lookup previously only returned a single
#    curve id; add the missing rows for Area1/Sub2 and Area2/Sub6 so each
#    org row maps to its own CurveID (1, 2, 3).
# ---------------------------------------------------------------------------
$wsLookup = $wb.Worksheets.Item("RespCurveLookup")
$wsLookup.Range("A3").Value = "Area1"
$wsLookup.Range("B3").Value = "Sub2"
$wsLookup.Range("C3").Value = "DM"
$wsLookup.Range("D3").Value = "Voice"
$wsLookup.Range("E3").Value = 2

$wsLookup.Range("A4").Value = "Area2"
$wsLookup.Range("B4").Value = "Sub6"
$wsLookup.Range("C4").Value = "DM"
$wsLookup.Range("D4").Value = "Voice"
$wsLookup.Range("E4").Value = 3

$wsLookup.Range("E5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Make "DatePeriods" the Instructions/first tab shown when the workbook
#    is opened (was "RespCurve" before), updating each sheet's selection to
#    match the new authored state.
# ---------------------------------------------------------------------------
$wsOrg.Range("C15").Select() | Out-Null

$wsRespCurve = $wb.Worksheets.Item("RespCurve")
$wsRespCurve.Range("E25").Select() | Out-Null

$wsDates.Range("C21").Select() | Out-Null
$wsDates.Activate() | Out-Null

$wb.Save()
